$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 11080.2
$ws.Range("I116").Value = 9901.5
$ws.Range("K116").Value = 9901.5
$ws.Range("M116").Value = -6459.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 50
$ws.Range("N4").Value = -282
$ws.Range("M4").ClearContents()

$ws.Range("H32").Value = 18870292
$ws.Range("I32").Value = 19232624
$ws.Range("K32").Value = 19232624
$ws.Range("M32").Value = -19232337

$ws.Range("H34").Value = 243799.8
$ws.Range("J34").Value = 262499.75
$ws.Range("L34").Value = 262499.75
$ws.Range("N34").Value = -263041.75

$ws.Range("H61").Value = 3604.75
$ws.Range("I61").Value = 3532.3484
$ws.Range("K61").Value = 3532.3484
$ws.Range("M61").Value = -3320.3484

$ws.Range("H74").Value = 13890399
$ws.Range("I74").Value = 13890399
$ws.Range("K74").Value = 13890399
$ws.Range("M74").Value = -13889525

$ws.Range("H77").Value = 13890399
$ws.Range("I77").Value = 13890399
$ws.Range("K77").Value = 69451995
$ws.Range("M77").Value = -69447627

$ws.Range("H112").Value = 59534.25
$ws.Range("J112").Value = 59534.25
$ws.Range("L112").Value = 59534.25
$ws.Range("N112").Value = -62488.25

$ws.Range("H136").Value = 3604.75
$ws.Range("I136").Value = 3532.3484
$ws.Range("K136").Value = 10597.0452
$ws.Range("M136").Value = -8047.0452

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 82500
$ws.Range("J109").Value = 82500
$ws.Range("L109").Value = 82500
$ws.Range("N109").Value = -85274

$ws.Range("H134").Value = 3134.375
$ws.Range("I134").Value = 1485.5385
$ws.Range("K134").Value = 4456.6155
$ws.Range("M134").Value = -1921.6155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33822.734
$ws.Range("I31").Value = 2228.4092
$ws.Range("K31").Value = 2228.4092
$ws.Range("M31").Value = -1933.4092

$ws.Range("H34").Value = 33822.734
$ws.Range("I34").Value = 2228.4092
$ws.Range("K34").Value = 2228.4092
$ws.Range("M34").Value = -2026.4092

$ws.Range("H86").Value = 8316.083000000001
$ws.Range("I86").Value = 3699.25
$ws.Range("J86").Value = 10624.5
$ws.Range("K86").Value = 3699.25
$ws.Range("L86").Value = 10624.5
$ws.Range("M86").Value = -2576.25
$ws.Range("N86").Value = -12870.5

$ws.Range("H89").Value = 8316.083000000001
$ws.Range("I89").Value = 3699.25
$ws.Range("J89").Value = 10624.5
$ws.Range("K89").Value = 18496.25
$ws.Range("L89").Value = 53122.5
$ws.Range("M89").Value = -12880.25
$ws.Range("N89").Value = -64354.5

$ws.Range("H107").Value = 1064.5294
$ws.Range("I107").Value = 987.7273
$ws.Range("J107").Value = 1205.3334
$ws.Range("K107").Value = 987.7273
$ws.Range("L107").Value = 1205.3334
$ws.Range("M107").Value = 932.2727
$ws.Range("N107").Value = -5045.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 13633.091
$ws.Range("J5").Value = 28203.2
$ws.Range("L5").Value = 84609.60000000001
$ws.Range("N5").Value = -84833.60000000001

$ws.Range("H11").Value = 674.5
$ws.Range("I11").Value = 674.5
$ws.Range("K11").Value = 2023.5
$ws.Range("M11").Value = -1883.5

$ws.Range("H129").Value = 5954940
$ws.Range("J129").Value = 7578657
$ws.Range("L129").Value = 22735971
$ws.Range("N129").Value = -22745971

$ws.Range("H131").Value = 8334758.5
$ws.Range("I131").Value = 8334638
$ws.Range("J131").Value = 8334939.5
$ws.Range("K131").Value = 25003914
$ws.Range("L131").Value = 25004818.5
$ws.Range("M131").Value = -24998874
$ws.Range("N131").Value = -25014898.5

$ws.Range("H132").Value = 3491.842
$ws.Range("I132").Value = 1434.6
$ws.Range("J132").Value = 5777.6665
$ws.Range("K132").Value = 12911.4
$ws.Range("L132").Value = 51998.9985
$ws.Range("M132").Value = -10381.4
$ws.Range("N132").Value = -57058.9985

$ws.Range("H135").Value = 13633.091
$ws.Range("J135").Value = 28203.2
$ws.Range("L135").Value = 253828.8
$ws.Range("N135").Value = -258898.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 563.2
$ws.Range("I2").Value = 87
$ws.Range("J2").Value = 1277.5
$ws.Range("K2").Value = 87
$ws.Range("L2").Value = 1277.5
$ws.Range("M2").Value = 26
$ws.Range("N2").Value = -1503.5

$ws.Range("H70").Value = 11208.546
$ws.Range("I70").Value = 8051.364
$ws.Range("J70").Value = 14365.728
$ws.Range("K70").Value = 8051.364
$ws.Range("L70").Value = 14365.728
$ws.Range("M70").Value = -7781.364
$ws.Range("N70").Value = -14905.728

$ws.Range("H73").Value = 11208.546
$ws.Range("I73").Value = 8051.364
$ws.Range("J73").Value = 14365.728
$ws.Range("K73").Value = 8051.364
$ws.Range("L73").Value = 14365.728
$ws.Range("M73").Value = -7115.364
$ws.Range("N73").Value = -16237.728

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 26000.715
$ws.Range("I20").Value = 25668.334
$ws.Range("K20").Value = 25668.334
$ws.Range("M20").Value = -25442.334

$ws.Range("H100").Value = 10828.59
$ws.Range("J100").Value = 12949.523
$ws.Range("L100").Value = 12949.523
$ws.Range("N100").Value = -14031.523

$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws.Range("H136").Value = 10100.0625
$ws.Range("I136").Value = 3732.5557
$ws.Range("J136").Value = 18286.857
$ws.Range("K136").Value = 11197.6671
$ws.Range("L136").Value = 54860.571
$ws.Range("M136").Value = -8647.667099999999
$ws.Range("N136").Value = -59960.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 28332.834
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 28332.834
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H132").Value = 12724.066
$ws.Range("I132").Value = 9693.714
$ws.Range("K132").Value = 29081.142
$ws.Range("M132").Value = -26551.142

$ws.Range("H136").Value = 1670.3784
$ws.Range("I136").Value = 1188.8611
$ws.Range("J136").Value = 19005
$ws.Range("K136").Value = 3566.5833
$ws.Range("L136").Value = 57015
$ws.Range("M136").Value = -1016.5833
$ws.Range("N136").Value = -62115
